$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.855.17'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '3.579.53'
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.73'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '189.38'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.633'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.00%  '
$ws.Range("D8").Value = '3.576.77'
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  -1.54%  '
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.83'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000303'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.65'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").Value = '4.153.52'
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.82'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.10%  '
$ws.Range("D17").Value = '3.572.38'
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("D18").Value = '69.779.13'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.67'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("E21").Value = '  -1.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '475.88'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.40'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.03'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.92%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.39'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.81%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.39'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.39%  '
$ws.Range("E27").Value = '  -3.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.05'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.97%  '
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("E30").Value = '  +5.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.42'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.09'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '582.65'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.11'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.13%  '
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("D38").Value = '0.0₃0799'
$ws.Range("E38").Value = '  -3.66%  '
$ws.Range("E39").Value = '  -4.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.23'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +18.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.90'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.19%  '
$ws.Range("E42").Value = '  -6.35%  '
$ws.Range("D43").Value = '3.244.77'
$ws.Range("E43").Value = '  -2.65%  '
$ws.Range("E44").Value = '  -6.08%  '
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("E46").Value = '  -1.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.37'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.44'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.83%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("E51").Value = '  -5.73%  '
